$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add time-of-day components to the date serials in column A for rows 211-267
# (row 249 already had no time component change per the source diff)
$ws.Range("A211").Value = 44138.416666666664
$ws.Range("A212").Value = 44138.5
$ws.Range("A213").Value = 44138.552083333336
$ws.Range("A214").Value = 44138.583333333336
$ws.Range("A215").Value = 44138.604166666664
$ws.Range("A216").Value = 44138.65625
$ws.Range("A217").Value = 44138.75
$ws.Range("A218").Value = 44138.791666666664
$ws.Range("A219").Value = 44138.798611111109
$ws.Range("A220").Value = 44138.805555555555
$ws.Range("A221").Value = 44138.809027777781
$ws.Range("A222").Value = 44138.8125
$ws.Range("A223").Value = 44138.815972222219
$ws.Range("A224").Value = 44138.822916666664
$ws.Range("A225").Value = 44138.826388888891
$ws.Range("A226").Value = 44138.833333333336
$ws.Range("A227").Value = 44138.84375
$ws.Range("A228").Value = 44138.850694444445
$ws.Range("A229").Value = 44138.857638888891
$ws.Range("A230").Value = 44138.864583333336
$ws.Range("A231").Value = 44138.868055555555
$ws.Range("A232").Value = 44138.875
$ws.Range("A233").Value = 44138.888888888891
$ws.Range("A234").Value = 44138.892361111109
$ws.Range("A235").Value = 44138.895833333336
$ws.Range("A236").Value = 44138.899305555555
$ws.Range("A237").Value = 44138.902777777781
$ws.Range("A238").Value = 44138.916666666664
$ws.Range("A239").Value = 44138.923611111109
$ws.Range("A240").Value = 44138.930555555555
$ws.Range("A241").Value = 44138.940972222219
$ws.Range("A242").Value = 44138.947916666664
$ws.Range("A243").Value = 44138.958333333336
$ws.Range("A244").Value = 44138.961805555555
$ws.Range("A245").Value = 44138.96875
$ws.Range("A246").Value = 44138.979166666664
$ws.Range("A247").Value = 44138.986111111109
$ws.Range("A248").Value = 44138.993055555555
$ws.Range("A250").Value = 44139.006944444445
$ws.Range("A251").Value = 44139.010416666664
$ws.Range("A252").Value = 44139.020833333336
$ws.Range("A253").Value = 44139.034722222219
$ws.Range("A254").Value = 44139.041666666664
$ws.Range("A255").Value = 44139.048611111109
$ws.Range("A256").Value = 44139.375
$ws.Range("A257").Value = 44139.416666666664
$ws.Range("A258").Value = 44139.458333333336
$ws.Range("A259").Value = 44139.510416666664
$ws.Range("A260").Value = 44139.5625
$ws.Range("A261").Value = 44139.565972222219
$ws.Range("A262").Value = 44139.638888888891
$ws.Range("A263").Value = 44139.708333333336
$ws.Range("A264").Value = 44140.017361111109
$ws.Range("A265").Value = 44140.430555555555
$ws.Range("A266").Value = 44140.652777777781
$ws.Range("A267").Value = 44141.388888888891

# Update selected cell to reflect where editing left off
$ws.Range("A257").Select()
